$wb = $excel.ActiveWorkbook

# Debit
$ws = $wb.Worksheets.Item("Debit")
$ws.Range("B2").Value = "Fri Nov 14 13:53:07 EST 2025"
$ws.Range("B3").Value = "Fri Nov 14 13:53:12 EST 2025"
$ws.Range("B4").Value = "Fri Nov 14 13:53:18 EST 2025"
$ws.Range("B5").Value = "Fri Nov 14 13:53:24 EST 2025"
$ws.Range("B6").Value = "Fri Nov 14 13:53:31 EST 2025"
$ws.Range("B7").Value = "Fri Nov 14 13:53:36 EST 2025"
$ws.Range("B8").Value = "Fri Nov 14 13:53:42 EST 2025"

# Debit-ZeroDollar
$ws = $wb.Worksheets.Item("Debit-ZeroDollar")
$ws.Range("B2").Value = "Fri Nov 14 13:53:48 EST 2025"
$ws.Range("B3").Value = "Fri Nov 14 13:53:54 EST 2025"
$ws.Range("B4").Value = "Fri Nov 14 13:54:00 EST 2025"
$ws.Range("B5").Value = "Fri Nov 14 13:54:06 EST 2025"
$ws.Range("B6").Value = "Fri Nov 14 13:54:11 EST 2025"
$ws.Range("B7").Value = "Fri Nov 14 13:54:17 EST 2025"
$ws.Range("B8").Value = "Fri Nov 14 13:54:24 EST 2025"

# Debit-Void
$ws = $wb.Worksheets.Item("Debit-Void")
$ws.Range("B2").Value = "Fri Nov 14 13:54:30 EST 2025"
$ws.Range("B3").Value = "Fri Nov 14 13:54:43 EST 2025"
$ws.Range("B4").Value = "Fri Nov 14 13:54:55 EST 2025"
$ws.Range("B5").Value = "Fri Nov 14 13:55:07 EST 2025"
$ws.Range("B6").Value = "Fri Nov 14 13:55:19 EST 2025"
$ws.Range("B7").Value = "Fri Nov 14 13:55:30 EST 2025"
$ws.Range("B8").Value = "Fri Nov 14 13:55:42 EST 2025"

# Debit-Credit
$ws = $wb.Worksheets.Item("Debit-Credit")
$ws.Range("B2").Value = "Fri Nov 14 13:55:54 EST 2025"
$ws.Range("B3").Value = "Fri Nov 14 13:56:07 EST 2025"
$ws.Range("B4").Value = "Fri Nov 14 13:56:17 EST 2025"
$ws.Range("B5").Value = "Fri Nov 14 13:56:29 EST 2025"
$ws.Range("B6").Value = "Fri Nov 14 13:56:41 EST 2025"
$ws.Range("B7").Value = "Fri Nov 14 13:56:52 EST 2025"
$ws.Range("B8").Value = "Fri Nov 14 13:57:04 EST 2025"

# Debit-Credit-Void
$ws = $wb.Worksheets.Item("Debit-Credit-Void")
$ws.Range("B2").Value = "Fri Nov 14 13:57:15 EST 2025"
$ws.Range("B3").Value = "Fri Nov 14 13:57:33 EST 2025"
$ws.Range("B4").Value = "Fri Nov 14 13:57:50 EST 2025"
$ws.Range("B5").Value = "Fri Nov 14 13:58:08 EST 2025"
$ws.Range("B6").Value = "Fri Nov 14 13:58:26 EST 2025"
$ws.Range("B7").Value = "Fri Nov 14 13:58:43 EST 2025"
$ws.Range("B8").Value = "Fri Nov 14 13:59:00 EST 2025"

# Debit-MRF
$ws = $wb.Worksheets.Item("Debit-MRF")
$ws.Range("B2").Value = "Sun Nov 09 12:11:40 EST 2025"
$ws.Range("B3").Value = "Sun Nov 09 12:11:46 EST 2025"
$ws.Range("B4").Value = "Sun Nov 09 12:11:53 EST 2025"
$ws.Range("B5").Value = "Sun Nov 09 12:11:59 EST 2025"
$ws.Range("B6").Value = "Sun Nov 09 12:12:05 EST 2025"
$ws.Range("B7").Value = "Sun Nov 09 12:12:12 EST 2025"
$ws.Range("B8").Value = "Sun Nov 09 12:12:18 EST 2025"
$ws.Range("B9").Value = "Sun Nov 09 12:12:25 EST 2025"
$ws.Range("B10").Value = "Sun Nov 09 12:12:31 EST 2025"
$ws.Range("B11").Value = "Sun Nov 09 12:12:38 EST 2025"
$ws.Range("B12").Value = "Sun Nov 09 12:12:45 EST 2025"
$ws.Range("B13").Value = "Sun Nov 09 12:12:51 EST 2025"
$ws.Range("B14").Value = "Sun Nov 09 12:12:58 EST 2025"
$ws.Range("B15").Value = "Sun Nov 09 12:13:04 EST 2025"
$ws.Range("B16").Value = "Sun Nov 09 12:13:10 EST 2025"
$ws.Range("B17").Value = "Sun Nov 09 12:13:17 EST 2025"
$ws.Range("B18").Value = "Sun Nov 09 12:13:24 EST 2025"
$ws.Range("B19").Value = "Sun Nov 09 12:13:31 EST 2025"
$ws.Range("B20").Value = "Sun Nov 09 12:13:37 EST 2025"

# Void-MRF
$ws = $wb.Worksheets.Item("Void-MRF")
$ws.Range("B2").Value = "Sun Nov 09 12:13:44 EST 2025"
$ws.Range("B3").Value = "Sun Nov 09 12:13:49 EST 2025"
$ws.Range("B4").Value = "Sun Nov 09 12:13:55 EST 2025"
$ws.Range("B5").Value = "Sun Nov 09 12:14:00 EST 2025"
$ws.Range("B6").Value = "Sun Nov 09 12:14:07 EST 2025"
$ws.Range("B7").Value = "Sun Nov 09 12:14:12 EST 2025"
$ws.Range("B8").Value = "Sun Nov 09 12:14:17 EST 2025"
$ws.Range("B9").Value = "Sun Nov 09 12:14:22 EST 2025"
$ws.Range("B10").Value = "Sun Nov 09 12:14:28 EST 2025"
$ws.Range("B11").Value = "Sun Nov 09 12:14:34 EST 2025"

# Credit-MRF
$ws = $wb.Worksheets.Item("Credit-MRF")
$ws.Range("B2").Value = "Sun Nov 09 12:14:40 EST 2025"
$ws.Range("B3").Value = "Sun Nov 09 12:14:45 EST 2025"
$ws.Range("B4").Value = "Sun Nov 09 12:14:52 EST 2025"
$ws.Range("B5").Value = "Sun Nov 09 12:14:57 EST 2025"
$ws.Range("B6").Value = "Sun Nov 09 12:15:03 EST 2025"
$ws.Range("B7").Value = "Sun Nov 09 12:15:09 EST 2025"
$ws.Range("B8").Value = "Sun Nov 09 12:15:15 EST 2025"
$ws.Range("B9").Value = "Sun Nov 09 12:15:21 EST 2025"
$ws.Range("B10").Value = "Sun Nov 09 12:15:26 EST 2025"
$ws.Range("B11").Value = "Sun Nov 09 12:15:31 EST 2025"
$ws.Range("B12").Value = "Sun Nov 09 12:15:36 EST 2025"
$ws.Range("B13").Value = "Sun Nov 09 12:15:41 EST 2025"

# Debit-RemID-Pipe
$ws = $wb.Worksheets.Item("Debit-RemID-Pipe")
$ws.Range("B2").Value = "Fri Nov 14 14:00:41 EST 2025"
$ws.Range("B3").Value = "Fri Nov 14 14:00:47 EST 2025"
$ws.Range("B4").Value = "Fri Nov 14 14:00:54 EST 2025"
$ws.Range("B5").Value = "Fri Nov 14 14:01:00 EST 2025"
$ws.Range("B6").Value = "Fri Nov 14 14:01:06 EST 2025"
$ws.Range("B7").Value = "Fri Nov 14 14:01:12 EST 2025"
$ws.Range("B8").Value = "Fri Nov 14 14:01:18 EST 2025"

# DebitVoid-RemID-Pipe
$ws = $wb.Worksheets.Item("DebitVoid-RemID-Pipe")
$ws.Range("B2").Value = "Fri Nov 14 14:01:25 EST 2025"
$ws.Range("B3").Value = "Fri Nov 14 14:01:37 EST 2025"
$ws.Range("B4").Value = "Fri Nov 14 14:01:49 EST 2025"
$ws.Range("B5").Value = "Fri Nov 14 14:02:00 EST 2025"
$ws.Range("B6").Value = "Fri Nov 14 14:02:14 EST 2025"
$ws.Range("B7").Value = "Fri Nov 14 14:02:26 EST 2025"
$ws.Range("B8").Value = "Fri Nov 14 14:02:37 EST 2025"

# DebitCredit-RemID-Pipe
$ws = $wb.Worksheets.Item("DebitCredit-RemID-Pipe")
$ws.Range("B2").Value = "Fri Nov 14 13:59:19 EST 2025"
$ws.Range("B3").Value = "Fri Nov 14 13:59:30 EST 2025"
$ws.Range("B4").Value = "Fri Nov 14 13:59:41 EST 2025"
$ws.Range("B5").Value = "Fri Nov 14 13:59:52 EST 2025"
$ws.Range("B6").Value = "Fri Nov 14 14:00:04 EST 2025"
$ws.Range("B7").Value = "Fri Nov 14 14:00:19 EST 2025"
$ws.Range("B8").Value = "Fri Nov 14 14:00:30 EST 2025"
